# Applies the "results with fixed workflow" update:
# For each sheet (NBR, BAR), the Cutoff/Reaction_number data (columns B and C)
# is shifted up by 4 rows (new row N gets the values that used to live at
# row N+4), column A keeps its 0-based running index, and the now-unused
# trailing 4 rows (17-20) are removed, shrinking the sheet from A1:C20 to
# A1:C16.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {

    # Capture the old B/C values from rows 6..20 before we start overwriting
    # rows 2..16 (row N new = row N+4 old).
    $oldB = @{}
    $oldC = @{}
    for ($r = 6; $r -le 20; $r++) {
        $oldB[$r] = $ws.Cells.Item($r, 2).Value2
        $oldC[$r] = $ws.Cells.Item($r, 3).Value2
    }

    # Write the shifted values into rows 2..16 (column A is left untouched).
    for ($newRow = 2; $newRow -le 16; $newRow++) {
        $srcRow = $newRow + 4
        $ws.Cells.Item($newRow, 2).Value = $oldB[$srcRow]
        $ws.Cells.Item($newRow, 3).Value = $oldC[$srcRow]
    }

    # Drop the now-redundant trailing rows 17-20 (deleting shifts dimension
    # down to A1:C16 automatically).
    $ws.Rows.Item(17).Resize(4).Delete()
}
